# ---------------------------------------------------------------------------
# Commit: "Thu, Jul 09, 2020  8:07:50 AM"
#
# Two independent changes:
#   1. The table on slide 5 switches from the deck's custom "Table_0" style
#      to the built-in PowerPoint table style {439AEEBE-784E-494C-B618-DEF2152B119A}.
#   2. The deck's theme (ppt/theme/theme1.xml, the slide master's theme,
#      originally the "Integral" / "Red Violet" color set) is recolored to
#      the stock Office theme palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------------
$slide = $p.Slides.Item(5)
$table = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $table = $candidate.Table
        break
    }
}
$table.ApplyStyle("{439AEEBE-784E-494C-B618-DEF2152B119A}")

# --- 2. Recolor the presentation theme to the Office palette ---------------
function ConvertTo-OleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme's fixed 12-slot layout:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$master      = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-OleColor($officeColors[$i - 1])
}
